$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.8
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 1.13
$ws.Range("L3").Value = 1.63
$ws.Range("U3").Value = 6.5
$ws.Range("AD3").Value = 9.5
$ws.Range("AE3").Value = 26
$ws.Range("AI3").Value = 81

# Row 4
$ws.Range("G4").Value = 2.15
$ws.Range("H4").Value = 2.8
$ws.Range("I4").Value = 4.1
$ws.Range("J4").Value = 1.13
$ws.Range("L4").Value = 1.63
$ws.Range("U4").Value = 8
$ws.Range("W4").Value = 19

# Row 5
$ws.Range("G5").Value = 2.65
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 2.5
$ws.Range("Q5").Value = 2.55
$ws.Range("R5").Value = 1.62
$ws.Range("S5").Value = 2.02
$ws.Range("T5").Value = 8.75
$ws.Range("U5").Value = 13.5
$ws.Range("V5").Value = 9.75
$ws.Range("W5").Value = 30
$ws.Range("X5").Value = 22
$ws.Range("Y5").Value = 29
$ws.Range("AB5").Value = 12.5
$ws.Range("AE5").Value = 13.5
$ws.Range("AF5").Value = 9.25
$ws.Range("AG5").Value = 28
$ws.Range("AH5").Value = 19.5
$ws.Range("AI5").Value = 26

# Row 8
$ws.Range("G8").Value = 3.35
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 2.02
$ws.Range("M8").Value = 3.35
$ws.Range("S8").Value = 2.1
$ws.Range("U8").Value = 20
$ws.Range("V8").Value = 11.5
$ws.Range("W8").Value = 50
$ws.Range("X8").Value = 28
$ws.Range("Y8").Value = 30
$ws.Range("AB8").Value = 12.5
$ws.Range("AC8").Value = 50
$ws.Range("AD8").Value = 8.5
$ws.Range("AE8").Value = 10.5
$ws.Range("AG8").Value = 19
$ws.Range("AH8").Value = 15

# Row 11
$ws.Range("K11").Value = 10

# Row 14
$ws.Range("Q14").Value = 2.8
$ws.Range("AE14").Value = 11
$ws.Range("AF14").Value = 9
